$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# 1) "A" -> "Ahmad" in the table header row (row 1, column 3): append "hmad"
#    after the existing "A" run's text so the cell reads "Ahmad".
$nameCell = $t.Cell(1, 3)
$nameCell.Range.InsertAfter("hmad")

# 2) Move the "_GoBack" bookmark from the trailing title paragraph to the
#    empty paragraph in row 2, column 3 of the table (the "A" column of the
#    "Setting up GitHub (10H)" row). Re-fetch the table/cell since the text
#    insertion above may have shifted cached ranges.
$t = $d.Tables.Item(1)
$targetCell = $t.Cell(2, 3)
$d.Bookmarks.Add("_GoBack", $targetCell.Range)
